$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A82").NumberFormat = "@"
$ws.Range("A82").Value = "11/22/2025"
$ws.Range("B82").Value = 0.2103874149996593
$ws.Range("C82").Value = 0.7896125850003407
$ws.Range("A82").ClearFormats()
